$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs, Gal, Galr2, ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gal"
$ws.Range("C2").Value = "Galr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.020062
$ws.Range("H2").Value = 3.060186
$ws.Range("I2").Value = 0.02710251114301613
$ws.Range("J2").Value = 0.02710251114301613
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9773006666666667
$ws.Range("N2").Value = 2.931902
$ws.Range("O2").Value = 0.2404260176543542
$ws.Range("P2").Value = 0.2404260176543543
$ws.Range("Q2").Value = 0.9969072726413334
$ws.Range("R2").Value = 8.972165453772002
$ws.Range("S2").Value = 0.006516148822548128
$ws.Range("T2").Value = 0.00651614882254813

# Row 3: ECs, Gal, Galr2, FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gal"
$ws.Range("C3").Value = "Galr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.020062
$ws.Range("H3").Value = 3.060186
$ws.Range("I3").Value = 0.02710251114301613
$ws.Range("J3").Value = 0.02710251114301613
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.102325333333333
$ws.Range("N3").Value = 6.306976000000001
$ws.Range("O3").Value = 0.5171936589700435
$ws.Range("P3").Value = 0.5171936589700435
$ws.Range("Q3").Value = 2.144502184170667
$ws.Range("R3").Value = 19.300519657536
$ws.Range("S3").Value = 0.01401724690533289
$ws.Range("T3").Value = 0.01401724690533289

# Row 4: ECs, Gal, Galr2, sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gal"
$ws.Range("C4").Value = "Galr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.020062
$ws.Range("H4").Value = 3.060186
$ws.Range("I4").Value = 0.02710251114301613
$ws.Range("J4").Value = 0.02710251114301613
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9852446666666667
$ws.Range("N4").Value = 2.955734
$ws.Range("O4").Value = 0.2423803233756023
$ws.Range("P4").Value = 0.2423803233756023
$ws.Range("Q4").Value = 1.005010645169333
$ws.Range("R4").Value = 9.045095806524001
$ws.Range("S4").Value = 0.006569115415135113
$ws.Range("T4").Value = 0.006569115415135115

# Row 5: sCs, Gal, Galr2, ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Gal"
$ws.Range("C5").Value = "Galr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 36.61711466666667
$ws.Range("H5").Value = 109.851344
$ws.Range("I5").Value = 0.9728974888569838
$ws.Range("J5").Value = 0.9728974888569839
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9773006666666667
$ws.Range("N5").Value = 2.931902
$ws.Range("O5").Value = 0.2404260176543542
$ws.Range("P5").Value = 0.2404260176543543
$ws.Range("Q5").Value = 35.78593057514311
$ws.Range("R5").Value = 322.073375176288
$ws.Range("S5").Value = 0.2339098688318061
$ws.Range("T5").Value = 0.2339098688318061

# Row 6: sCs, Gal, Galr2, FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gal"
$ws.Range("C6").Value = "Galr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 36.61711466666667
$ws.Range("H6").Value = 109.851344
$ws.Range("I6").Value = 0.9728974888569838
$ws.Range("J6").Value = 0.9728974888569839
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.102325333333333
$ws.Range("N6").Value = 6.306976000000001
$ws.Range("O6").Value = 0.5171936589700435
$ws.Range("P6").Value = 0.5171936589700435
$ws.Range("Q6").Value = 76.98108779730489
$ws.Range("R6").Value = 692.8297901757441
$ws.Range("S6").Value = 0.5031764120647106
$ws.Range("T6").Value = 0.5031764120647106

# Row 7: sCs, Gal, Galr2, sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gal"
$ws.Range("C7").Value = "Galr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 36.61711466666667
$ws.Range("H7").Value = 109.851344
$ws.Range("I7").Value = 0.9728974888569838
$ws.Range("J7").Value = 0.9728974888569839
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9852446666666667
$ws.Range("N7").Value = 2.955734
$ws.Range("O7").Value = 0.2423803233756023
$ws.Range("P7").Value = 0.2423803233756023
$ws.Range("Q7").Value = 36.07681693405511
$ws.Range("R7").Value = 324.691352406496
$ws.Range("S7").Value = 0.2358112079604671
$ws.Range("T7").Value = 0.2358112079604672
